# Apply the edits described by the diff:
#  - Rename worksheet "loginTest" -> "ValidateLogin"
#  - Make "ValidateLogin" the active sheet/tab (it was "openAccount" before)
#  - Update the selected cell on "ValidateLogin" from C28 to D7
#  - "openAccount" keeps its own selection (B3) but is no longer the active tab

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("loginTest")
$loginSheet.Name = "ValidateLogin"

$loginSheet.Activate()
$loginSheet.Range("D7").Select()

$wb.Save()
